# Prok1-Prokr1.xlsx — refresh with new TPM numbers.
#
# The sending-cluster set changed from a single "MuSCs" sender (5 target
# rows, one of them "Inflammatory-Mac") to two senders — "ECs" and
# "MuSCs" — each crossed with 4 targets ("ECs","FAPs","MuSCs",
# "Resolving-Mac"); "Inflammatory-Mac" no longer appears anywhere.
# Ligand/receptor stay "Prok1"/"Prokr1" throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T, one tuple per output row (rows 2..9 on the sheet).
$rows = @(
    ("ECs","Prok1","Prokr1","ECs",1,0.3333333333333333,0.1903723333333333,0.571117,0.7691179579859432,0.7691179579859432,3,1,2.419247,7.257740999999999,0.3378166463347285,0.3378166463347285,0.4605576962996666,4.145019266696999,0.2598208492026259,0.2598208492026259),
    ("ECs","Prok1","Prokr1","FAPs",1,0.3333333333333333,0.1903723333333333,0.571117,0.7691179579859432,0.7691179579859432,3,1,4.158501,12.475503,0.5806810390173532,0.5806810390173534,0.7916635385390001,7.124971846850999,0.4466122149701825,0.4466122149701826),
    ("ECs","Prok1","Prokr1","MuSCs",1,0.3333333333333333,0.1903723333333333,0.571117,0.7691179579859432,0.7691179579859432,1,0.3333333333333333,0.02378566666666667,0.071357,0.00332136162374866,0.00332136162374866,0.004528132863222223,0.040753195769,0.002554518869790445,0.002554518869790446),
    ("ECs","Prok1","Prokr1","Resolving-Mac",1,0.3333333333333333,0.1903723333333333,0.571117,0.7691179579859432,0.7691179579859432,3,1,0.5598866666666668,1.67966,0.07818095302416966,0.07818095302416966,0.1065869311355556,0.95928238022,0.06013037494334432,0.06013037494334432),
    ("MuSCs","Prok1","Prokr1","ECs",3,1,0.057148,0.171444,0.2308820420140567,0.2308820420140567,3,1,2.419247,7.257740999999999,0.3378166463347285,0.3378166463347285,0.138255127556,1.244296148004,0.07799579713210253,0.07799579713210254),
    ("MuSCs","Prok1","Prokr1","FAPs",3,1,0.057148,0.171444,0.2308820420140567,0.2308820420140567,3,1,4.158501,12.475503,0.5806810390173532,0.5806810390173534,0.237650015148,2.138850136332,0.1340688240471707,0.1340688240471707),
    ("MuSCs","Prok1","Prokr1","MuSCs",3,1,0.057148,0.171444,0.2308820420140567,0.2308820420140567,1,0.3333333333333333,0.02378566666666667,0.071357,0.00332136162374866,0.00332136162374866,0.001359303278666667,0.012233729508,0.0007668427539582137,0.0007668427539582139),
    ("MuSCs","Prok1","Prokr1","Resolving-Mac",3,1,0.057148,0.171444,0.2308820420140567,0.2308820420140567,3,1,0.5598866666666668,1.67966,0.07818095302416966,0.07818095302416966,0.03199640322666667,0.28796762904,0.01805057808082534,0.01805057808082534)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowValues = $rows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
